$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spanish month abbreviations (with trailing period) used in column C (Mes)
$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

# Replace the numeric month values in column C (rows 6-85) with their
# textual abbreviation equivalents.
for ($r = 6; $r -le 85; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $monthNum = [int]$cell.Value2
    $cell.Value = $monthNames[$monthNum]
}

# Update the "ND No disponible" legend note to include a trailing period.
$ws.Cells.Item(88, 2).Value = "ND No disponible."
